# 037 Week 21/22 Update
# Fills in the WK22 scores that were missing for a number of players on the
# Sunday Pairs sheet and the Thursday Singles sheet, and updates the two
# manually-entered handicap values on the HANDICAPS sheet. All dependent
# totals (AB/AC weekly+aggregate sums, the Z running totals on THURSDAY
# SINGLES and the COUNTIF "played" counts on the hidden summary sheet) are
# plain formulas, so they recalculate automatically once the inputs below
# are written.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1 (WINTER BEST PAIRS COMP) - WK22 column is "Y"
# ---------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")

$sheet1.Range("Y13").Value = 36
# This particular score is highlighted in red in the source workbook.
$sheet1.Range("Y13").Font.Color = 255

$sheet1.Range("Y17").Value = 30
$sheet1.Range("Y18").Value = 31
$sheet1.Range("Y19").Value = 28
$sheet1.Range("Y20").Value = 32
$sheet1.Range("Y21").Value = 35
$sheet1.Range("Y27").Value = 25
$sheet1.Range("Y28").Value = 35
$sheet1.Range("Y31").Value = 33

# ---------------------------------------------------------------------
# THURSDAY SINGLES - WK21 column is "V"
# ---------------------------------------------------------------------
$sheet2 = $wb.Worksheets.Item("THURSDAY SINGLES")

$sheet2.Range("V7").Value = 35
# This particular score is highlighted in red in the source workbook.
$sheet2.Range("V7").Font.Color = 255

$sheet2.Range("V10").Value = 30
$sheet2.Range("V12").Value = 27
$sheet2.Range("V16").Value = 27

# ---------------------------------------------------------------------
# HANDICAPS - BAZ MASON's handicap dropped by one shot on both the
# Sunday (pairs) and Thursday (singles) competitions.
# ---------------------------------------------------------------------
$sheet3 = $wb.Worksheets.Item("HANDICAPS")

$sheet3.Range("B2").Value = 11
$sheet3.Range("C2").Value = 11
